$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Columns A-F: literal text values for new session rows 1034-1046 ---
$ws.Range("A1034").Value = '2024-03-19'
$ws.Range("B1034").Value = '16:30'
$ws.Range("C1034").Value = '19:15'
$ws.Range("D1034").Value = '2h 45m'
$ws.Range("E1034").Value = '#python'
$ws.Range("F1034").Value = 'nwtimetrackingmanager v3.0.0'

$ws.Range("A1035").Value = '2024-03-19'
$ws.Range("B1035").Value = '20:30'
$ws.Range("C1035").Value = '21:15'
$ws.Range("D1035").Value = '0h 45m'
$ws.Range("E1035").Value = '#python'
$ws.Range("F1035").Value = 'nwtimetrackingmanager v3.0.0'

$ws.Range("A1036").Value = '2024-03-21'
$ws.Range("B1036").Value = '08:00'
$ws.Range("C1036").Value = '08:45'
$ws.Range("D1036").Value = '0h 45m'
$ws.Range("E1036").Value = '#python'
$ws.Range("F1036").Value = 'nwreadinglistmanager v3.0.0'

$ws.Range("A1037").Value = '2024-03-21'
$ws.Range("B1037").Value = '17:15'
$ws.Range("C1037").Value = '17:45'
$ws.Range("D1037").Value = '0h 30m'
$ws.Range("E1037").Value = '#python'
$ws.Range("F1037").Value = 'nwreadinglistmanager v3.0.0'

$ws.Range("A1038").Value = '2024-03-22'
$ws.Range("B1038").Value = '08:00'
$ws.Range("C1038").Value = '08:45'
$ws.Range("D1038").Value = '0h 45m'
$ws.Range("E1038").Value = '#python'
$ws.Range("F1038").Value = 'nwreadinglistmanager v3.0.0'

$ws.Range("A1039").Value = '2024-03-22'
$ws.Range("B1039").Value = '17:00'
$ws.Range("C1039").Value = '17:45'
$ws.Range("D1039").Value = '0h 45m'
$ws.Range("E1039").Value = '#python'
$ws.Range("F1039").Value = 'nwreadinglistmanager v3.0.0'

$ws.Range("A1040").Value = '2024-03-22'
$ws.Range("B1040").Value = '23:30'
$ws.Range("C1040").Value = '00:45'
$ws.Range("D1040").Value = '1h 15m'
$ws.Range("E1040").Value = '#python'
$ws.Range("F1040").Value = 'nwreadinglistmanager v3.0.0'

$ws.Range("A1041").Value = '2024-03-23'
$ws.Range("B1041").Value = '10:30'
$ws.Range("C1041").Value = '11:30'
$ws.Range("D1041").Value = '1h 00m'
$ws.Range("E1041").Value = '#python'
$ws.Range("F1041").Value = 'nwreadinglistmanager v3.0.0'

$ws.Range("A1042").Value = '2024-03-23'
$ws.Range("B1042").Value = '16:00'
$ws.Range("C1042").Value = '20:30'
$ws.Range("D1042").Value = '4h 30m'
$ws.Range("E1042").Value = '#python'
$ws.Range("F1042").Value = 'nwreadinglistmanager v3.0.0'

$ws.Range("A1043").Value = '2024-03-23'
$ws.Range("B1043").Value = '23:00'
$ws.Range("C1043").Value = '23:30'
$ws.Range("D1043").Value = '0h 30m'
$ws.Range("E1043").Value = '#python'
$ws.Range("F1043").Value = 'nwreadinglistmanager v3.0.0'

$ws.Range("A1044").Value = '2024-03-24'
$ws.Range("B1044").Value = '11:30'
$ws.Range("C1044").Value = '13:45'
$ws.Range("D1044").Value = '2h 15m'
$ws.Range("E1044").Value = '#python'
$ws.Range("F1044").Value = 'nwreadinglistmanager v3.1.0'

$ws.Range("A1045").Value = '2024-03-24'
$ws.Range("B1045").Value = '14:30'
$ws.Range("C1045").Value = '17:00'
$ws.Range("D1045").Value = '2h 30m'
$ws.Range("E1045").Value = '#python'
$ws.Range("F1045").Value = 'nwreadinglistmanager v3.1.0'

$ws.Range("A1046").Value = '2024-03-28'
$ws.Range("B1046").Value = '11:00'
$ws.Range("C1046").Value = '11:30'
$ws.Range("D1046").Value = '0h 30m'
$ws.Range("E1046").Value = '#python'
$ws.Range("F1046").Value = 'nwreadinglistmanager v3.1.0'

# --- Columns G/H: True/False flags, copied from existing text-typed cells so they
#     stay literal text ("True"/"False") like the rest of the column instead of
#     being coerced into native booleans ---
$trueSrc = $ws.Range("G48")
$falseSrc = $ws.Range("G2")

for ($r = 1034; $r -le 1046; $r++) {
  $trueSrc.Copy($ws.Cells.Item($r, 7))
}

for ($r = 1034; $r -le 1035; $r++) {
  $trueSrc.Copy($ws.Cells.Item($r, 8))
}
for ($r = 1036; $r -le 1040; $r++) {
  $falseSrc.Copy($ws.Cells.Item($r, 8))
}
for ($r = 1041; $r -le 1043; $r++) {
  $trueSrc.Copy($ws.Cells.Item($r, 8))
}
for ($r = 1044; $r -le 1045; $r++) {
  $falseSrc.Copy($ws.Cells.Item($r, 8))
}
$trueSrc.Copy($ws.Cells.Item(1046, 8))

# --- Columns I/J: YEAR/MONTH helper formulas, entered range-at-a-time (matches
#     the shared-formula grouping Excel itself produces on a multi-cell fill) ---
$ws.Range("I1034:I1035").Formula = "=YEAR(A1034)"
$ws.Range("J1034:J1035").Formula = "=MONTH(A1034)"
$ws.Range("I1036:I1037").Formula = "=YEAR(A1036)"
$ws.Range("J1036:J1037").Formula = "=MONTH(A1036)"
$ws.Range("I1038:I1039").Formula = "=YEAR(A1038)"
$ws.Range("J1038:J1039").Formula = "=MONTH(A1038)"
$ws.Range("I1040:I1042").Formula = "=YEAR(A1040)"
$ws.Range("J1040:J1042").Formula = "=MONTH(A1040)"
$ws.Range("I1043:I1045").Formula = "=YEAR(A1043)"
$ws.Range("J1043:J1045").Formula = "=MONTH(A1043)"
$ws.Range("I1046").Formula = "=YEAR(A1046)"
$ws.Range("J1046").Formula = "=MONTH(A1046)"

# --- Extend the blank trailing rows: 1047/1048 already existed blank; add
#     1049-1063 (same blank style) by copying the last existing blank row down ---
$ws.Range("A1048:J1048").Copy($ws.Range("A1049:J1063"))

# --- Selection state: land on A1047, matching the post-edit cursor position ---
$ws.Range("A1047").Select()
